$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "list" data validations that were present on columns A-E
$ws.Cells.Validation.Delete()

# Column C ("Birth") values were stored as text dates; convert them to real
# Excel date/time serial values and apply a date-time number format so they
# export as numeric dates instead of shared strings.
$c5 = $ws.Range("C5")
$c5.NumberFormat = "yyyy-MM-dd HH:mm"
$c5.Value = (Get-Date -Year 2021 -Month 1 -Day 1 -Hour 8 -Minute 32 -Second 0)

$c6 = $ws.Range("C6")
$c6.NumberFormat = "yyyy-MM-dd HH:mm"
$c6.Value = (Get-Date -Year 2022 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0)
